$wb = $excel.ActiveWorkbook

# Add the new worksheet "3.25-4.1" after the existing "3.17-3.24" sheet
$oldSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $oldSheet)
$newSheet.Name = "3.25-4.1"

# Column widths
$newSheet.Columns.Item(1).ColumnWidth = 27.875
$newSheet.Columns.Item(2).ColumnWidth = 71.125

# Fill data
$newSheet.Range("A1").Value = "Yangwei"
$newSheet.Range("B1").Value = "Forge System"
$newSheet.Range("B2").Value = "Customizable Weapon And booster storage"
$newSheet.Range("B3").Value = "Avatar weapon replacement"
$newSheet.Range("B4").Value = "Item Enhancement System"
$newSheet.Range("B5").Value = "Xp & Coins collection"
$newSheet.Range("B6").Value = "Items Buy"
$newSheet.Range("B7").Value = "LevelUp (Configs and logic)"

$newSheet.Range("A9").Value = "JiaNan"
$newSheet.Range("B10").Value = "Avatar animation complete (start hit, gournd hit, air hit)"
$newSheet.Range("B11").Value = "Floor enemy, Fly enemy design"

$newSheet.Range("A13").Value = "Langyefan"
$newSheet.Range("B13").Value = "Items(Heros, Weapons, Boosters, Skills) Config"
$newSheet.Range("B14").Value = "Global constants definition"

# Set the active selection on the new sheet
$newSheet.Range("B14").Select()

# Update selection on the first sheet to span A1:B10
$oldSheet.Range("A1:B10").Select()

# Make the new sheet active / selected tab
$newSheet.Activate()
